$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("G2").Value = 16.58130233333334
$ws.Range("H2").Value = 49.74390700000001
$ws.Range("I2").Value = 0.6205214785234227
$ws.Range("J2").Value = 0.6205214785234225
$ws.Range("M2").Value = 46.63275166666667
$ws.Range("N2").Value = 139.898255
$ws.Range("O2").Value = 0.9158911059585902
$ws.Range("P2").Value = 0.9158911059585902
$ws.Range("Q2").Value = 773.2317540202541
$ws.Range("R2").Value = 6959.085786182286
$ws.Range("S2").Value = 0.5683301032358772
$ws.Range("T2").Value = 0.5683301032358771

# Row 3
$ws.Range("G3").Value = 16.58130233333334
$ws.Range("H3").Value = 49.74390700000001
$ws.Range("I3").Value = 0.6205214785234227
$ws.Range("J3").Value = 0.6205214785234225
$ws.Range("O3").Value = 0.05441917700612491
$ws.Range("P3").Value = 0.05441917700612491
$ws.Range("Q3").Value = 45.94283688861056
$ws.Range("R3").Value = 413.485531997495
$ws.Range("S3").Value = 0.03376826817586847
$ws.Range("T3").Value = 0.03376826817586847

# Row 4
$ws.Range("G4").Value = 16.58130233333334
$ws.Range("H4").Value = 49.74390700000001
$ws.Range("I4").Value = 0.6205214785234227
$ws.Range("J4").Value = 0.6205214785234225
$ws.Range("M4").Value = 0.849605
$ws.Range("N4").Value = 2.548815
$ws.Range("O4").Value = 0.01668667696558362
$ws.Range("P4").Value = 0.01668667696558362
$ws.Range("Q4").Value = 14.08755736891167
$ws.Range("R4").Value = 126.788016320205
$ws.Range("S4").Value = 0.01035444146232669
$ws.Range("T4").Value = 0.01035444146232669

# Row 5
$ws.Range("G5").Value = 16.58130233333334
$ws.Range("H5").Value = 49.74390700000001
$ws.Range("I5").Value = 0.6205214785234227
$ws.Range("J5").Value = 0.6205214785234225
$ws.Range("M5").Value = 0.662052
$ws.Range("N5").Value = 1.986156
$ws.Range("O5").Value = 0.01300304006970129
$ws.Range("P5").Value = 0.0130030400697013
$ws.Range("Q5").Value = 10.977684372388
$ws.Range("R5").Value = 98.79915935149202
$ws.Range("S5").Value = 0.008068665649350356
$ws.Range("T5").Value = 0.008068665649350356

# Row 6
$ws.Range("I6").Value = 0.02090549052511678
$ws.Range("J6").Value = 0.02090549052511678
$ws.Range("M6").Value = 46.63275166666667
$ws.Range("N6").Value = 139.898255
$ws.Range("O6").Value = 0.9158911059585902
$ws.Range("P6").Value = 0.9158911059585902
$ws.Range("Q6").Value = 26.05032970954555
$ws.Range("R6").Value = 234.45296738591
$ws.Range("S6").Value = 0.01914715283765604
$ws.Range("T6").Value = 0.01914715283765604

# Row 7
$ws.Range("I7").Value = 0.02090549052511678
$ws.Range("J7").Value = 0.02090549052511678
$ws.Range("O7").Value = 0.05441917700612491
$ws.Range("P7").Value = 0.05441917700612491
$ws.Range("S7").Value = 0.001137659589286197
$ws.Range("T7").Value = 0.001137659589286197

# Row 8
$ws.Range("I8").Value = 0.02090549052511678
$ws.Range("J8").Value = 0.02090549052511678
$ws.Range("M8").Value = 0.849605
$ws.Range("N8").Value = 2.548815
$ws.Range("O8").Value = 0.01668667696558362
$ws.Range("P8").Value = 0.01668667696558362
$ws.Range("Q8").Value = 0.4746125755366666
$ws.Range("R8").Value = 4.27151317983
$ws.Range("S8").Value = 0.0003488431671996929
$ws.Range("T8").Value = 0.0003488431671996929

# Row 9
$ws.Range("I9").Value = 0.02090549052511678
$ws.Range("J9").Value = 0.02090549052511678
$ws.Range("M9").Value = 0.662052
$ws.Range("N9").Value = 1.986156
$ws.Range("O9").Value = 0.01300304006970129
$ws.Range("P9").Value = 0.0130030400697013
$ws.Range("Q9").Value = 0.3698403432879999
$ws.Range("R9").Value = 3.328563089592
$ws.Range("S9").Value = 0.0002718349309748542
$ws.Range("T9").Value = 0.0002718349309748543

# Row 10
$ws.Range("G10").Value = 9.581631
$ws.Range("H10").Value = 28.744893
$ws.Range("I10").Value = 0.3585730309514606
$ws.Range("J10").Value = 0.3585730309514606
$ws.Range("M10").Value = 46.63275166666667
$ws.Range("N10").Value = 139.898255
$ws.Range("O10").Value = 0.9158911059585902
$ws.Range("P10").Value = 0.9158911059585902
$ws.Range("Q10").Value = 446.817818984635
$ws.Range("R10").Value = 4021.360370861715
$ws.Range("S10").Value = 0.3284138498850571
$ws.Range("T10").Value = 0.3284138498850571

# Row 11
$ws.Range("G11").Value = 9.581631
$ws.Range("H11").Value = 28.744893
$ws.Range("I11").Value = 0.3585730309514606
$ws.Range("J11").Value = 0.3585730309514606
$ws.Range("O11").Value = 0.05441917700612491
$ws.Range("P11").Value = 0.05441917700612491
$ws.Range("Q11").Value = 26.548415878945
$ws.Range("R11").Value = 238.935742910505
$ws.Range("S11").Value = 0.01951324924097024
$ws.Range("T11").Value = 0.01951324924097024

# Row 12
$ws.Range("G12").Value = 9.581631
$ws.Range("H12").Value = 28.744893
$ws.Range("I12").Value = 0.3585730309514606
$ws.Range("J12").Value = 0.3585730309514606
$ws.Range("M12").Value = 0.849605
$ws.Range("N12").Value = 2.548815
$ws.Range("O12").Value = 0.01668667696558362
$ws.Range("P12").Value = 0.01668667696558362
$ws.Range("Q12").Value = 8.140601605755
$ws.Range("R12").Value = 73.26541445179501
$ws.Range("S12").Value = 0.005983392336057242
$ws.Range("T12").Value = 0.005983392336057242

# Row 13
$ws.Range("G13").Value = 9.581631
$ws.Range("H13").Value = 28.744893
$ws.Range("I13").Value = 0.3585730309514606
$ws.Range("J13").Value = 0.3585730309514606
$ws.Range("M13").Value = 0.662052
$ws.Range("N13").Value = 1.986156
$ws.Range("O13").Value = 0.01300304006970129
$ws.Range("P13").Value = 0.0130030400697013
$ws.Range("Q13").Value = 6.343537966812
$ws.Range("R13").Value = 57.09184170130801
$ws.Range("S13").Value = 0.004662539489376084
$ws.Range("T13").Value = 0.004662539489376085

# Remove rows 14-17 (Resolving-Mac as sending cluster) entirely
$ws.Range("A14:T17").Delete()

Write-Host "Edit complete"